$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# This section's first-page header/footer differ from the default
# (non-first-page) header/footer, so there are three picture-bearing
# header/footer "stories" in the document:
#   Headers.Item(2) (first page)  -> BTec logo    (was named image1.jpg)
#   Footers.Item(1) (default)     -> Pearson logo (was named image2.png)
#   Footers.Item(2) (first page)  -> Pearson logo (was named image2.png)
#
# The three embedded logo pictures get renumbered: the two Pearson-logo
# pictures flip from "image2.png" to "image1.png", and the BTec-logo
# picture flips from "image1.jpg" to "image2.jpg".
#
# Renaming is routed through Selection.InlineShapes (select the picture,
# then rename the shape at the selection) rather than renaming the
# HeaderFooter.Range.InlineShapes item directly in place, which keeps each
# rename targeting a freshly-resolved object.

function Rename-Logo($story, $altTextMatch, $newName) {
    for ($i = 1; $i -le $story.Range.InlineShapes.Count; $i++) {
        $candidate = $story.Range.InlineShapes.Item($i)
        if ($candidate.AlternativeText -like $altTextMatch) {
            [void]$candidate.Select()
            $word.Selection.InlineShapes.Item(1).Name = $newName
        }
    }
}

# First-page header: BTec_Logo-Orange -> image1.jpg becomes image2.jpg
Rename-Logo $sec.Headers.Item(2) "BTec_Logo-Orange" "image2.jpg"

# Default (primary) footer: PearsonLogo.png -> image2.png becomes image1.png
Rename-Logo $sec.Footers.Item(1) "*PearsonLogo.png" "image1.png"

# First-page footer: PearsonLogo.png -> image2.png becomes image1.png
Rename-Logo $sec.Footers.Item(2) "*PearsonLogo.png" "image1.png"

Write-Output "Renamed header/footer logo pictures."
